# REVER_DailyTracker_MONISHA.xlsx - "Add files via upload"
#
# The DEC-2020 sheet (already the active sheet / active tab) gets three new
# daily entries (rows 7, 8, 9 -> worksheet rows 8, 9, 10) filled in under the
# existing task log, and one now-unused blank spacer row removed further
# down the sheet (shifting the legend block at the bottom up by one row).
# The current selection also moves from C13 to D13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the three new log rows -----------------------------------
# Copy the existing row-4 formatting (borders/alignment/number formats for
# columns A:F, and the plain fill for G) down onto the new rows so the new
# cells pick up the same styles already used for the rest of the table.
$ws.Range("A4:F4").Copy()
$ws.Range("A8:F10").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G4").Copy()
$ws.Range("G8").PasteSpecial(-4122)       # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 44172
$ws.Range("C8").Value = " nMVAR"
$ws.Range("D8").Value = " nMVAR issues checking"
$ws.Range("E8").Value = 0.5
$ws.Range("F8").Value = "WIP"
$ws.Rows("8").RowHeight = 25.95

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 44173
$ws.Range("C9").Value = " nMVAR"
$ws.Range("D9").Value = " nMVAR issues checking"
$ws.Range("E9").Value = 0.75
$ws.Range("F9").Value = "WIP"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 44174
$ws.Range("C10").Value = " nMVAR"
$ws.Range("D10").Value = " nMVAR issues checking"
$ws.Range("E10").Value = 0.9
$ws.Range("F10").Value = "WIP"

# --- Remove the now-redundant blank spacer row ------------------------
# This shifts the legend block (Status/WIP/Pending/Completed/Hold) that used
# to start at row 19 up to row 18, and shrinks the used range from G23 to
# G22.
$ws.Rows("16").Delete()

# --- Match the saved cursor position -----------------------------------
$ws.Range("D13").Select() | Out-Null
